$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.225.03"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "1.657.82"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'218.42"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'0.5222"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.2673"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "'0.06336"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").Value = "'21.18"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").Value = "'0.07734"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "'4.443"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "1.652.29"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "1.881.68"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "'0.5485"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").Value = "0.0₅8275"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "'65.00"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "26.233.81"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "'4.672"
$ws.Range("D21").Value = "'193.31"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").Value = "'10.17"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").Value = "'6.121"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'138.00"
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("D26").Value = "'0.1245"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").Value = "'7.252"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("D28").Value = "'16.18"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'1.429"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").Value = "'0.06033"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "'1.284"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'3.567"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "'3.343"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").Value = "'0.9840"
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "'2.770"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "'0.5937"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("D39").Value = "'0.01597"
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("D40").Value = "'5.966"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "'0.8643"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.003"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.042.05"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "'99.74"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'57.28"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "'0.05177"
$ws.Range("D51").Value = "'1.470"
$ws.Range("E51").Value = "  +3.78%  "
